$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update check-in / check-out dates for Kansas City (row 4) and Seattle (row 2)
$ws.Range("B4").Value = " 11/15/2021"
$ws.Range("B2").Value = " 02/15/2022"
$ws.Range("C4").Value = " 01/05/2022"
$ws.Range("C2").Value = " 02/19/2022"

# Update the active selection to B2
$ws.Range("B2").Select()
